$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the "OECD - Total" row (row 40) entirely, shifting rows below up.
$ws.Rows.Item(40).Delete()

# Update the view: scroll so A16 is the top-left visible cell, and select D43
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("D43").Select()
